$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 0: stage the date value (B8) so we can paste it as TEXT (not an
# auto-converted date serial) into the new B13/C13 cells later, even after
# the original rows have been deleted. PasteSpecial-values preserves the
# shared-string "t=s" type and the original cell style.
$ws.Range("B8").Copy()
$ws.Range("Z1").PasteSpecial(-4163)   # xlPasteValues

# --- Phase 1: capture all other source text we still need from the region
# that is about to be restructured (rows 13-27), while it's still intact.
$txt_A15 = $ws.Range("A15").Value2    # "Programa resumido:"
$txt_A16 = $ws.Range("A16").Value2    # "Short syllabus:"
$txt_A17 = $ws.Range("A17").Value2    # "Programa:"
$txt_B13 = $ws.Range("B13").Value2    # "6495737 - Durval Rodrigues Junior"
$txt_A18 = $ws.Range("A18").Value2    # "Syllabus:"
$txt_A19 = $ws.Range("A19").Value2    # "Avaliação:"
$txt_A20 = $ws.Range("A20").Value2    # "Método:"
$txt_B14 = $ws.Range("B14").Value2    # "1643715 - Paulo Atsushi Suzuki"
$txt_A21 = $ws.Range("A21").Value2    # "Critério:"
$txt_B20 = $ws.Range("B20").Value2    # "Experimentos desenvolvidos..."
$txt_A22 = $ws.Range("A22").Value2    # "Norma de recuperação:"
$txt_B21 = $ws.Range("B21").Value2    # "Média aritmética..."
$txt_A23 = $ws.Range("A23").Value2    # "Bibliografia:"
$txt_B22 = $ws.Range("B22").Value2    # "Aplicação de uma prova..."
$txt_A24 = $ws.Range("A24").Value2    # "Requisitos:"
$txt_B25 = $ws.Range("B25").Value2    # "LOB1021 ..."
$txt_B26 = $ws.Range("B26").Value2    # "LOM3016 ..."
$txt_B27 = $ws.Range("B27").Value2    # "LOM3246 ..."

# --- Phase 2: remove the old rows 13-27 entirely.
$ws.Range("A13:C27").EntireRow.Delete()

# --- Phase 3: rebuild rows 13-25 with the new layout.

# row 13: Programa resumido: | 01/01/2012 | 01/01/2012
$ws.Range("A13").Value = $txt_A15
$ws.Range("Z1").Copy()
$ws.Range("B13:C13").PasteSpecial(-4163)
$ws.Rows.Item(13).RowHeight = 60

# row 14: Short syllabus:
$ws.Range("A14").Value = $txt_A16
$ws.Rows.Item(14).RowHeight = 60

# row 15: Programa: | 6495737 - Durval Rodrigues Junior | (same)
$ws.Range("A15").Value = $txt_A17
$ws.Range("B15").Value = $txt_B13
$ws.Range("C15").Value = $txt_B13
$ws.Rows.Item(15).RowHeight = 120

# row 16: Syllabus:
$ws.Range("A16").Value = $txt_A18
$ws.Rows.Item(16).RowHeight = 120

# row 17: Avaliação:
$ws.Range("A17").Value = $txt_A19

# row 18: Método: | 1643715 - Paulo Atsushi Suzuki | (same)
$ws.Range("A18").Value = $txt_A20
$ws.Range("B18").Value = $txt_B14
$ws.Range("C18").Value = $txt_B14
$ws.Rows.Item(18).RowHeight = 60

# row 19: Critério: | Experimentos desenvolvidos... | (same)
$ws.Range("A19").Value = $txt_A21
$ws.Range("B19").Value = $txt_B20
$ws.Range("C19").Value = $txt_B20
$ws.Rows.Item(19).RowHeight = 60

# row 20: Norma de recuperação: | Média aritmética... | (same)
$ws.Range("A20").Value = $txt_A22
$ws.Range("B20").Value = $txt_B21
$ws.Range("C20").Value = $txt_B21
$ws.Rows.Item(20).RowHeight = 60

# row 21: Bibliografia: | Aplicação de uma prova... | (same)
$ws.Range("A21").Value = $txt_A23
$ws.Range("B21").Value = $txt_B22
$ws.Range("C21").Value = $txt_B22
$ws.Rows.Item(21).RowHeight = 120

# row 22: Requisitos:
$ws.Range("A22").Value = $txt_A24

# row 23: LOB1021 ...
$ws.Range("B23").Value = $txt_B25
$ws.Range("C23").Value = $txt_B25
$ws.Rows.Item(23).RowHeight = 30

# row 24: LOM3016 ...
$ws.Range("B24").Value = $txt_B26
$ws.Range("C24").Value = $txt_B26
$ws.Rows.Item(24).RowHeight = 30

# row 25: LOM3246 ...
$ws.Range("B25").Value = $txt_B27
$ws.Range("C25").Value = $txt_B27
$ws.Rows.Item(25).RowHeight = 30

# --- Phase 4: drop the staging cell so it doesn't bleed into the used range.
$ws.Range("Z1").Clear()

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
